# Fruta / hortaliza, semanal
# Insert two new weekly observation rows at row 72 (pushing the existing
# rows 72-206 down to 74-208), then populate the two new rows with the
# latest "Ajo" price data for Macroferia Regional de Talca.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 72. Excel shifts rows
# 72..206 down to 74..208 (and duplicates formatting/styles from the
# row above, which matches the existing date-column style).
$ws.Rows.Item(72).Insert()
$ws.Rows.Item(72).Insert()

# New row 72
$ws.Range("A72").Value = 5
$ws.Range("B72").Value = "Macroferia Regional de Talca"
$ws.Range("C72").Value = "Maule"
$ws.Range("D72").Value = 44533
$ws.Range("E72").Value = 7
$ws.Range("F72").Value = 100112003
$ws.Range("G72").Value = "Ajo"
$ws.Range("H72").Value = "Chino"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 200
$ws.Range("K72").Value = 19000
$ws.Range("L72").Value = 19000
$ws.Range("M72").Value = 19000
$ws.Range("N72").Value = "`$/caja 10 kilos"
$ws.Range("O72").Value = "China"
$ws.Range("P72").Value = 1900
$ws.Range("Q72").Value = 10
$ws.Range("R72").Value = "Hortaliza"

# New row 73
$ws.Range("A73").Value = 5
$ws.Range("B73").Value = "Macroferia Regional de Talca"
$ws.Range("C73").Value = "Maule"
$ws.Range("D73").Value = 44533
$ws.Range("E73").Value = 7
$ws.Range("F73").Value = 100112003
$ws.Range("G73").Value = "Ajo"
$ws.Range("H73").Value = "Chino"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 200
$ws.Range("K73").Value = 20000
$ws.Range("L73").Value = 20000
$ws.Range("M73").Value = 20000
$ws.Range("N73").Value = "`$/malla 10 kilos"
$ws.Range("O73").Value = "China"
$ws.Range("P73").Value = 2000
$ws.Range("Q73").Value = 10
$ws.Range("R73").Value = "Hortaliza"
